# Update cryptos list values (Price / Volume(1h) / Coin / Link) to reflect refreshed market data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.624.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.12%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.789.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.16%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.91"
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.790.96"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.523"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "

# Row 10
$ws.Range("E10").Value = "  -2.21%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.39%  "

# Row 12
$ws.Range("E12").Value = "  -1.97%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000266"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.24%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.422.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.25%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.786.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.96%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.602.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.31%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.36%  "

# Row 19
$ws.Range("E19").Value = "  -2.89%  "

# Row 20
$ws.Range("E20").Value = "  -0.29%  "

# Row 21
$ws.Range("E21").Value = "  +2.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "472.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.707"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.53%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000145"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.89%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.34%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.79%  "

# Row 29
$ws.Range("E29").Value = "  +0.13%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.934.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.23%  "

# Row 31
$ws.Range("E31").Value = "  -4.21%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.68%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.88%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.742.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.40%  "

# Row 38
$ws.Range("E38").Value = "  -2.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.63%  "

# Row 40
$ws.Range("E40").Value = "  +1.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.85"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.13%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.311"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.41%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.72%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.31%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "405.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.81%  "

# Row 49
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "41.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.06%  "

# Row 50
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.75%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.78%  "
